# Pac Man death sequence with sound done.
#
# The first TODO bullet ("Add player and ghost collision detection.")
# is replaced by the text that used to be the third bullet
# ("Have level restart. ..."), and the now-duplicate second/third
# bullets ("Add death state ..." and "Have level restart ...") are
# removed, leaving a single bullet behind (with its bookmark intact).

$d = $word.ActiveDocument

# 1. Swap the text of the first bullet for the "Have level restart..." text.
$d.Content.Find.Execute(
    "Add player and ghost collision detection.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Have level restart. Check with video or play game to see how level resets.",
    2) | Out-Null

# 2. Remove the two now-redundant bullet paragraphs that followed it
#    ("Add death state and have Pac-Man animate death." and the old
#    "Have level restart..." bullet), paragraph mark included.
$first = $d.Paragraphs.Item(4)
$last  = $d.Paragraphs.Item(5)
$d.Range($first.Range.Start, $last.Range.End).Delete() | Out-Null
